$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The weekly consolidation now has two additional daily records for this
# market/variety ("Sandia" @ Terminal Hortofrutícola Agro Chillán). Insert
# two new rows right before the current row 48 and push the existing data
# (rows 48-158) down to rows 50-160.
$ws.Rows("48:49").Insert()

# New row 48: Extra quality, fecha 44592
$ws.Cells.Item(48, 1).Value  = 7
$ws.Cells.Item(48, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(48, 3).Value  = "Ñuble"
$ws.Cells.Item(48, 4).Value  = 44592
$ws.Cells.Item(48, 5).Value  = 16
$ws.Cells.Item(48, 6).Value  = 100112028
$ws.Cells.Item(48, 7).Value  = "Sandia"
$ws.Cells.Item(48, 8).Value  = "Sin especificar"
$ws.Cells.Item(48, 9).Value  = "Extra"
$ws.Cells.Item(48, 10).Value = 300
$ws.Cells.Item(48, 11).Value = 2000
$ws.Cells.Item(48, 12).Value = 2000
$ws.Cells.Item(48, 13).Value = 2000
$ws.Cells.Item(48, 14).Value = "$/unidad"
$ws.Cells.Item(48, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(48, 16).Value = 2000
$ws.Cells.Item(48, 17).Value = 1
$ws.Cells.Item(48, 18).Value = "Hortaliza"

# New row 49: Primera quality, fecha 44592
$ws.Cells.Item(49, 1).Value  = 7
$ws.Cells.Item(49, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(49, 3).Value  = "Ñuble"
$ws.Cells.Item(49, 4).Value  = 44592
$ws.Cells.Item(49, 5).Value  = 16
$ws.Cells.Item(49, 6).Value  = 100112028
$ws.Cells.Item(49, 7).Value  = "Sandia"
$ws.Cells.Item(49, 8).Value  = "Sin especificar"
$ws.Cells.Item(49, 9).Value  = "Primera"
$ws.Cells.Item(49, 10).Value = 300
$ws.Cells.Item(49, 11).Value = 1500
$ws.Cells.Item(49, 12).Value = 1800
$ws.Cells.Item(49, 13).Value = 1650
$ws.Cells.Item(49, 14).Value = "$/unidad"
$ws.Cells.Item(49, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(49, 16).Value = 1650
$ws.Cells.Item(49, 17).Value = 1
$ws.Cells.Item(49, 18).Value = "Hortaliza"
